$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AV3").Value = -73.70999999999999
$ws.Range("AW3").Value = 376.29
$ws.Range("AS4").Value = 63.93
$ws.Range("AW4").Value = 394.71
$ws.Range("AS5").Value = 90.08
$ws.Range("AW5").Value = 413.1300000000001
$ws.Range("AS6").Value = 92.66
$ws.Range("AV6").Value = 92.11
$ws.Range("AW6").Value = 505.2400000000001
$ws.Range("AS7").Value = 92.48
$ws.Range("AV7").Value = 92.11
$ws.Range("AW7").Value = 597.35
$ws.Range("AS8").Value = 45.5
$ws.Range("AV8").Value = -73.70999999999999
$ws.Range("AW8").Value = 523.64
$ws.Range("AS9").Value = 94.69
$ws.Range("AV9").Value = 92.11
$ws.Range("AW9").Value = 615.75
$ws.Range("AS10").Value = 94.51000000000001
$ws.Range("AV10").Value = 92.11
$ws.Range("AW10").Value = 707.86
$ws.Range("AV11").Value = -73.70999999999999
$ws.Range("AW11").Value = 634.15
$ws.Range("AS12").Value = 93.40000000000001
$ws.Range("AV12").Value = 92.11
$ws.Range("AW12").Value = 726.26
$ws.Range("AS13").Value = 42.74
$ws.Range("AW13").Value = 744.6799999999999
$ws.Range("AS14").Value = 87.69
$ws.Range("AV14").Value = -73.70999999999999
$ws.Range("AW14").Value = 670.9699999999999
$ws.Range("AS15").Value = 93.03
$ws.Range("AV15").Value = 92.11
$ws.Range("AW15").Value = 763.0799999999999
$ws.Range("AS16").Value = 93.22
$ws.Range("AV16").Value = 92.11
$ws.Range("AW16").Value = 855.1899999999999
